$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "43.180.60"
$ws.Cells.Item(2, 5).Value = "  +0.46%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.321.20"
$ws.Cells.Item(3, 5).Value = "  +0.98%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.02%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "303.27"
$ws.Cells.Item(5, 5).Value = "  +0.42%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "99.54"
$ws.Cells.Item(6, 5).Value = "  +0.59%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "0.518"
$ws.Cells.Item(9, 5).Value = "  +2.29%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "36.06"
$ws.Cells.Item(10, 5).Value = "  +5.64%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -0.57%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -0.86%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "17.74"
$ws.Cells.Item(13, 5).Value = "  -0.36%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  +2.26%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "2.682.59"
$ws.Cells.Item(15, 5).Value = "  +1.04%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "2.323.95"
$ws.Cells.Item(16, 5).Value = "  +0.11%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  -1.29%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "43.099.32"
$ws.Cells.Item(18, 5).Value = "  +0.52%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "12.92"
$ws.Cells.Item(19, 5).Value = "  +4.69%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  +2.86%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  +1.19%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "68.15"
$ws.Cells.Item(22, 5).Value = "  +0.60%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "240.70"
$ws.Cells.Item(23, 5).Value = "  +1.89%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -1.45%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +0.11%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "0.998"
$ws.Cells.Item(26, 5).Value = "  -0.18%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +3.27%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "168.38"
$ws.Cells.Item(28, 5).Value = "  +0.24%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "34.30"
$ws.Cells.Item(29, 5).Value = "  +1.96%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +0.54%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  -10.73%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  +2.45%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  -0.05%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "4.82"
$ws.Cells.Item(34, 5).Value = "  +5.66%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "17.70"
$ws.Cells.Item(35, 5).Value = "  +4.88%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -0.55%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "0.0698"
$ws.Cells.Item(37, 5).Value = "  +0.96%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "1.82"
$ws.Cells.Item(38, 5).Value = "  +2.42%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +0.39%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -0.45%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +0.43%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "1.995.84"
$ws.Cells.Item(42, 5).Value = "  +0.14%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  +1.68%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -4.97%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "10.09"
$ws.Cells.Item(45, 5).Value = "  +0.86%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "17.65"
$ws.Cells.Item(46, 5).Value = "  +0.62%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "2.85"
$ws.Cells.Item(47, 5).Value = "  +0.32%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "76.83"
$ws.Cells.Item(48, 5).Value = "  +9.95%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "55.08"
$ws.Cells.Item(49, 5).Value = "  -2.91%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "2.548.06"
$ws.Cells.Item(50, 5).Value = "  +0.80%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "Stacks"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(51, 4).Value = "1.54"
$ws.Cells.Item(51, 5).Value = "  +1.52%  "
